$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting existing rows 61-90 down to 62-91.
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new record's data.
$ws.Cells.Item(61, 1).Value = 7
$ws.Cells.Item(61, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(61, 3).Value = "Ñuble"
$ws.Cells.Item(61, 4).Value = 44960
$ws.Cells.Item(61, 4).NumberFormat = $ws.Cells.Item(62, 4).NumberFormat
$ws.Cells.Item(61, 5).Value = 16
$ws.Cells.Item(61, 6).Value = 100112037
$ws.Cells.Item(61, 7).Value = "Cebollín"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 200
$ws.Cells.Item(61, 11).Value = 600
$ws.Cells.Item(61, 12).Value = 600
$ws.Cells.Item(61, 13).Value = 600
$ws.Cells.Item(61, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(61, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(61, 16).Value = 100
$ws.Cells.Item(61, 17).Value = 6
$ws.Cells.Item(61, 18).Value = "Hortaliza"
